# Case and Fatality Demographics Data Updated
# Update the raw "Number" (column B) counts on each of the six sheets;
# the "%" column formulas (and the "Total"/"Grand Total" SUM cells)
# recalculate automatically.

$wb = $excel.ActiveWorkbook

# --- Cases by Age Group ---------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Age Group")
$ws.Range("B2").Value2 = 255
$ws.Range("B3").Value2 = 1234
$ws.Range("B4").Value2 = 3294
$ws.Range("B5").Value2 = 14024
$ws.Range("B6").Value2 = 15348
$ws.Range("B7").Value2 = 13492
$ws.Range("B8").Value2 = 11467
$ws.Range("B9").Value2 = 4163
$ws.Range("B10").Value2 = 2805
$ws.Range("B11").Value2 = 1669
$ws.Range("B12").Value2 = 1082
$ws.Range("B13").Value2 = 1679

# --- Cases by Gender -------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Gender")
$ws.Range("B2").Value2 = 24069
$ws.Range("B3").Value2 = 45567
$ws.Range("B4").Value2 = 890

# --- Cases by RaceEthnicity --------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws.Range("B2").Value2 = 920
$ws.Range("B3").Value2 = 11750
$ws.Range("B4").Value2 = 26927
$ws.Range("B5").Value2 = 365
$ws.Range("B6").Value2 = 22453
$ws.Range("B7").Value2 = 8111

# --- Fatalities by Age Group ------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B4").Value2 = 25
$ws.Range("B5").Value2 = 175
$ws.Range("B6").Value2 = 572
$ws.Range("B7").Value2 = 1704
$ws.Range("B8").Value2 = 3889
$ws.Range("B9").Value2 = 3260
$ws.Range("B10").Value2 = 4196
$ws.Range("B11").Value2 = 4746
$ws.Range("B12").Value2 = 4761
$ws.Range("B13").Value2 = 12647

# --- Fatalities by Gender ----------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value2 = 15209
$ws.Range("B3").Value2 = 20778

# --- Fatalities by Race-Ethnicity --------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value2 = 690
$ws.Range("B3").Value2 = 3413
$ws.Range("B4").Value2 = 16964
$ws.Range("B5").Value2 = 191
$ws.Range("B6").Value2 = 14709

# --- Restore each sheet's remembered selection (author re-saved after ------
# clicking around) and leave "Cases by Age Group" as the active tab. -------
$wb.Worksheets.Item("Cases by Gender").Range("E15").Select()
$wb.Worksheets.Item("Cases by RaceEthnicity").Range("C23").Select()
$wb.Worksheets.Item("Fatalities by Age Group").Range("H10").Select()
$wb.Worksheets.Item("Fatalities by Gender").Range("D15").Select()
$wb.Worksheets.Item("Fatalities by Race-Ethnicity").Range("B11").Select()
$wb.Worksheets.Item("Cases by Age Group").Range("A24").Select()
